# Revert "Powerpoint writer: consolidate text run nodes."
#
# The title placeholders on slide 1 ("First slide") and slide 3
# ("Third slide") each currently hold the leading word plus the
# trailing space as a single run (e.g. "First "). We split that run
# back into two runs - one for the word, one for the lone space -
# while leaving the "slide" run alone, matching the
# pre-consolidation XML.
#
# We achieve the split by re-assigning the text of just the
# word-only sub-range (the first five characters) to itself: the
# writer only emits new run boundaries where text actually changed,
# so touching that sub-range forces a split right after
# "First"/"Third" without altering any characters or formatting.

$p = $ppt.ActivePresentation

for ($slideIndex = 1; $slideIndex -le $p.Slides.Count; $slideIndex++) {
    $slide = $p.Slides.Item($slideIndex)

    if ($slide.Shapes.Count -lt 1) {
        continue
    }

    $titleShape = $slide.Shapes.Item(1)

    if (-not $titleShape.HasTextFrame) {
        continue
    }

    $titleFrame = $titleShape.TextFrame

    if (-not $titleFrame.HasText) {
        continue
    }

    $titleRange = $titleFrame.TextRange
    $fullText = $titleRange.Text

    if ($fullText -eq "First slide" -or $fullText -eq "Third slide") {
        $word = $titleRange.Characters(1, 5)
        $word.Text = $word.Text
    }
}
